$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 66, shifting rows 66-88 down to 67-89
$ws.Rows.Item(66).Insert()

# Populate the newly inserted row 66 with the new weekly data point
$ws.Range("A66").Value = 3
$ws.Range("B66").Value = "Femacal de La Calera"
$ws.Range("C66").Value = "Coquimbo"
$ws.Range("D66").Value = 44508
$ws.Range("E66").Value = 5
$ws.Range("F66").Value = 100112026
$ws.Range("G66").Value = "Haba"
$ws.Range("H66").Value = "Sin especificar"
$ws.Range("I66").Value = "Primera"
$ws.Range("J66").Value = 90
$ws.Range("K66").Value = 7000
$ws.Range("L66").Value = 7500
$ws.Range("M66").Value = 7278
$ws.Range("N66").Value = "$/malla 25 kilos"
$ws.Range("O66").Value = "Provincia de Quillota"
$ws.Range("P66").Value = 291
$ws.Range("Q66").Value = 25
$ws.Range("R66").Value = "Hortaliza"
